$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.644.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.449.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.110'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.347'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000173'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.900.21'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.458.20'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.435.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("E18").Value = '  -5.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '642.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.574.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0951'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.11%  '
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.54%  '
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '150.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.87%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.604'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0903'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.74%  '
